$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 560
$ws1.Range("F4").Value = 323
$ws1.Range("F5").Value = 401
$ws1.Range("F7").Value = 2363
$ws1.Range("F8").Value = 395
$ws1.Range("F9").Value = 6011
$ws1.Range("F11").Value = 389
$ws1.Range("F12").Value = 17

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 560
$ws4.Range("F5").Value = 323
$ws4.Range("F6").Value = 401
$ws4.Range("F10").Value = 2363
$ws4.Range("F11").Value = 395
$ws4.Range("F12").Value = 6011
$ws4.Range("F14").Value = 389
$ws4.Range("F16").Value = 17
